$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.950.25'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '1.814.91'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.19'
$ws.Range('E5').Value = '  -2.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5900'
$ws.Range('E6').Value = '  -3.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.004'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2744'
$ws.Range('E8').Value = '  -2.74%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06747'
$ws.Range('E9').Value = '  -4.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.88'
$ws.Range('E10').Value = '  -4.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07489'
$ws.Range('E11').Value = '  -2.08%  '
$ws.Range('D12').Value = '1.811.70'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.666'
$ws.Range('E13').Value = '  -3.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6219'
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009362'
$ws.Range('E15').Value = '  -6.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '74.58'
$ws.Range('E16').Value = '  -6.33%  '
$ws.Range('D17').Value = '28.737.93'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.421'
$ws.Range('E18').Value = '  -9.32%  '
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '207.68'
$ws.Range('E20').Value = '  -9.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.37'
$ws.Range('E21').Value = '  -3.92%  '
$ws.Range('E22').Value = '  -3.96%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '154.64'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1267'
$ws.Range('E25').Value = '  -2.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.762'
$ws.Range('E26').Value = '  -4.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.24'
$ws.Range('E27').Value = '  -3.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06356'
$ws.Range('E28').Value = '  -6.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.399'
$ws.Range('E29').Value = '  -5.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.428'
$ws.Range('E30').Value = '  -2.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.716'
$ws.Range('E31').Value = '  -3.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.670'
$ws.Range('E32').Value = '  -4.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.680'
$ws.Range('E33').Value = '  -3.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.045'
$ws.Range('E34').Value = '  -7.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.523'
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6299'
$ws.Range('E36').Value = '  -4.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.732'
$ws.Range('E37').Value = '  -1.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.408'
$ws.Range('E38').Value = '  -2.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01691'
$ws.Range('E39').Value = '  -4.43%  '
$ws.Range('D40').Value = '1.132.51'
$ws.Range('E40').Value = '  -8.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8672'
$ws.Range('E41').Value = '  -6.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.004'
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').Value = '1.972.11'
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.76'
$ws.Range('E44').Value = '  -1.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.20'
$ws.Range('E45').Value = '  -5.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000111'
$ws.Range('E46').Value = '  -4.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.571'
$ws.Range('E47').Value = '  -3.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05466'
$ws.Range('E48').Value = '  -1.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4514'
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.234'
$ws.Range('E50').Value = '  -4.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').Value = '  +0.02%  '
